# Re-order the "Recorded By" (column G) comma-separated list on the
# "Session Analysis Results" sheet: the last name/email in the list is
# moved to the front of the list (a left-rotation by one applied from
# the tail), e.g. "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Rows in the "Recorded By" (G) column whose value needs to be rotated.
$rows = @(2,3,6,7,10,11,12,13,14,15,17,18,19,20,21,22,24,26,28,29,32,33,36,37,38,39,40,41,43,44,45,46,47,48,50,52,54,55,58,59,62,63,64,65,66,67,69,70,71,72,73,74,76,78,83,84,85,86,90,92,93,94,96,99,101,109,110,111,112,116,118,119,120,122,125,127,135,136,137,138,142,144,145,146,148,151,153)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 7)
    $current = [string]$cell.Value2

    # Split on comma, trim whitespace around each entry.
    $parts = $current -split ','
    for ($i = 0; $i -lt $parts.Length; $i++) {
        $parts[$i] = $parts[$i].Trim()
    }

    if ($parts.Length -gt 1) {
        # Move the last entry to the front of the list.
        $last = $parts[$parts.Length - 1]
        $rest = $parts[0..($parts.Length - 2)]
        $newParts = @($last) + $rest
        $newValue = [string]::Join(', ', $newParts)
        $cell.Value = $newValue
    }
}
